$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates -------------------------------------------------
# "Volume 31   Number  6" -> "Volume 31   Number  7"
$ws.Range("A8").Value = "Volume 31   Number  7"
# "Report Covering the Week  2/5/2024  Through  2/11/2024"
#   -> "Report Covering the Week  2/12/2024  Through  2/18/2024"
$ws.Range("C9").Value = "Report Covering the Week  2/12/2024  Through  2/18/2024"

# --- Cells that flip between a numeric value and the "no data" text ------
# placeholders ("0" = shared text for nil count, "***.*" = shared text for
# an undefined percentage). Copy from a same-styled donor cell so the
# number format / style (s="14" text vs s="15"/"16" numeric) matches,
# then the donor's own value already equals the value we need.
$ws.Range("C14").Copy($ws.Range("C15"))   # -> text "0" (was numeric 1)
$ws.Range("F15").Copy($ws.Range("C22"))   # -> numeric 1   (was text "0")
$ws.Range("C14").Copy($ws.Range("D22"))   # -> text "0" (was numeric 1)
$ws.Range("E14").Copy($ws.Range("E22"))   # -> text "***.*" (was numeric -100)
$ws.Range("C14").Copy($ws.Range("C23"))   # -> text "0" (was numeric 1)
$ws.Range("C14").Copy($ws.Range("C26"))   # -> text "0" (was numeric 1)
$ws.Range("C14").Copy($ws.Range("F30"))   # -> text "0" (was numeric 1)

# --- Updated weekly crime-statistics figures ------------------------------
$ws.Range("L15").Value = -66.666666666666
$ws.Range("M15").Value = -75
$ws.Range("N15").Value = -91.666666666666
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = 66.666666666666
$ws.Range("F16").Value = 12
$ws.Range("H16").Value = -36.842105263157
$ws.Range("I16").Value = 30
$ws.Range("J16").Value = 34
$ws.Range("K16").Value = -11.764705882352
$ws.Range("L16").Value = -18.918918918918
$ws.Range("M16").Value = -23.076923076923
$ws.Range("N16").Value = -86.238532110091
$ws.Range("C17").Value = 4
$ws.Range("E17").Value = -20
$ws.Range("F17").Value = 18
$ws.Range("G17").Value = 15
$ws.Range("H17").Value = 20
$ws.Range("I17").Value = 36
$ws.Range("J17").Value = 37
$ws.Range("K17").Value = -2.702702702702
$ws.Range("L17").Value = 50
$ws.Range("M17").Value = 50
$ws.Range("N17").Value = -59.090909090909
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = -66.666666666666
$ws.Range("I18").Value = 14
$ws.Range("J18").Value = 19
$ws.Range("K18").Value = -26.315789473684
$ws.Range("L18").Value = -39.130434782608
$ws.Range("M18").Value = -60
$ws.Range("N18").Value = -94.696969696969
$ws.Range("C19").Value = 16
$ws.Range("E19").Value = 100
$ws.Range("F19").Value = 57
$ws.Range("G19").Value = 34
$ws.Range("H19").Value = 67.647058823529
$ws.Range("I19").Value = 85
$ws.Range("J19").Value = 59
$ws.Range("K19").Value = 44.067796610169
$ws.Range("L19").Value = 7.594936708860
$ws.Range("M19").Value = 123.684210526316
$ws.Range("N19").Value = -45.161290322580
$ws.Range("C20").Value = 9
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = 125
$ws.Range("F20").Value = 22
$ws.Range("H20").Value = 15.789473684210
$ws.Range("I20").Value = 33
$ws.Range("J20").Value = 32
$ws.Range("K20").Value = 3.125
$ws.Range("L20").Value = -47.619047619047
$ws.Range("M20").Value = 230
$ws.Range("N20").Value = -88
$ws.Range("C21").Value = 35
$ws.Range("D21").Value = 23
$ws.Range("E21").Value = 52.173913043478
$ws.Range("F21").Value = 119
$ws.Range("G21").Value = 101
$ws.Range("H21").Value = 17.821782178217
$ws.Range("I21").Value = 199
$ws.Range("J21").Value = 182
$ws.Range("K21").Value = 9.340659340659
$ws.Range("L21").Value = -13.100436681222
$ws.Range("M21").Value = 32.666666666666
$ws.Range("N21").Value = -80.604288499025
$ws.Range("F22").Value = 3
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 7
$ws.Range("K22").Value = 40
$ws.Range("L22").Value = -12.5
$ws.Range("M22").Value = 40
$ws.Range("E23").Value = -100
$ws.Range("G23").Value = 2
$ws.Range("H23").Value = 50
$ws.Range("J23").Value = 5
$ws.Range("K23").Value = -20
$ws.Range("M23").Value = 0
$ws.Range("C24").Value = 18
$ws.Range("D24").Value = 19
$ws.Range("E24").Value = -5.263157894736
$ws.Range("F24").Value = 96
$ws.Range("G24").Value = 81
$ws.Range("H24").Value = 18.518518518518
$ws.Range("I24").Value = 163
$ws.Range("J24").Value = 172
$ws.Range("K24").Value = -5.232558139534
$ws.Range("L24").Value = -6.321839080459
$ws.Range("M24").Value = 126.388888888889
$ws.Range("C25").Value = 15
$ws.Range("D25").Value = 7
$ws.Range("E25").Value = 114.285714285714
$ws.Range("F25").Value = 45
$ws.Range("G25").Value = 34
$ws.Range("H25").Value = 32.352941176470
$ws.Range("I25").Value = 70
$ws.Range("J25").Value = 51
$ws.Range("K25").Value = 37.254901960784
$ws.Range("L25").Value = 42.857142857142
$ws.Range("M25").Value = 7.692307692307
$ws.Range("L26").Value = -66.666666666666
$ws.Range("F27").Value = 3
$ws.Range("H27").Value = -40
$ws.Range("J27").Value = 8
$ws.Range("K27").Value = -25
$ws.Range("L27").Value = 0
$ws.Range("F28").Value = 1
$ws.Range("H28").Value = 0
$ws.Range("N28").Value = -85.294117647058
$ws.Range("F29").Value = 1
$ws.Range("H29").Value = 0
$ws.Range("N29").Value = -84.848484848484
